$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the data range stays formatted as Text so numeric-looking strings
# (e.g. "63.538.72", "0.591") are preserved exactly as text, matching the
# original inline-string cell values.
$ws.Range("B2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '63.538.72'
$ws.Range("E2").Value = '  -0.71%  '

# Row 3
$ws.Range("D3").Value = '2.719.29'
$ws.Range("E3").Value = '  -1.27%  '

# Row 4
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
$ws.Range("D5").Value = '561.67'
$ws.Range("E5").Value = '  -2.40%  '

# Row 6
$ws.Range("D6").Value = '159.48'
$ws.Range("E6").Value = '  +0.64%  '

# Row 7
$ws.Range("E7").Value = '  +0.13%  '

# Row 8
$ws.Range("D8").Value = '0.591'
$ws.Range("E8").Value = '  -2.15%  '

# Row 9
$ws.Range("D9").Value = '0.107'
$ws.Range("E9").Value = '  -2.39%  '

# Row 10
$ws.Range("E10").Value = '  +3.05%  '

# Row 11
$ws.Range("E11").Value = '  -0.45%  '

# Row 12
$ws.Range("D12").Value = '0.372'
$ws.Range("E12").Value = '  -2.98%  '

# Row 13
$ws.Range("D13").Value = '3.207.17'
$ws.Range("E13").Value = '  -1.18%  '

# Row 14
$ws.Range("D14").Value = '26.51'
$ws.Range("E14").Value = '  -1.41%  '

# Row 15
$ws.Range("D15").Value = '63.397.48'
$ws.Range("E15").Value = '  -0.38%  '

# Row 16
$ws.Range("D16").Value = '0.0000147'
$ws.Range("E16").Value = '  -2.83%  '

# Row 17
$ws.Range("D17").Value = '2.727.58'
$ws.Range("E17").Value = '  -1.22%  '

# Row 18
$ws.Range("D18").Value = '12.17'
$ws.Range("E18").Value = '  -0.06%  '

# Row 19
$ws.Range("D19").Value = '4.67'
$ws.Range("E19").Value = '  -3.56%  '

# Row 20
$ws.Range("D20").Value = '351.64'
$ws.Range("E20").Value = '  -1.35%  '

# Row 21
$ws.Range("D21").Value = '6.46'
$ws.Range("E21").Value = '  -3.99%  '

# Row 22
$ws.Range("D22").Value = '0.998'

# Row 23
$ws.Range("D23").Value = '0.509'
$ws.Range("E23").Value = '  -4.59%  '

# Row 24
$ws.Range("D24").Value = '63.84'
$ws.Range("E24").Value = '  -2.37%  '

# Row 25
$ws.Range("E25").Value = '  -1.22%  '

# Row 26
$ws.Range("E26").Value = '  +0.06%  '

# Row 27
$ws.Range("D27").Value = '8.22'
$ws.Range("E27").Value = '  -5.15%  '

# Row 28
$ws.Range("D28").Value = '0.0₃0892'
$ws.Range("E28").Value = '  -1.50%  '

# Row 29
$ws.Range("E29").Value = '  -0.05%  '

# Row 30
$ws.Range("D30").Value = '1.35'
$ws.Range("E30").Value = '  +6.63%  '

# Row 31
$ws.Range("D31").Value = '7.13'
$ws.Range("E31").Value = '  -1.67%  '

# Row 32
$ws.Range("D32").Value = '166.77'
$ws.Range("E32").Value = '  -2.32%  '

# Row 33
$ws.Range("D33").Value = '1.49'
$ws.Range("E33").Value = '  -0.05%  '

# Row 34
$ws.Range("D34").Value = '19.86'
$ws.Range("E34").Value = '  -1.68%  '

# Row 35
$ws.Range("B35").Value = 'USDe'
$ws.Range("C35").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.02%  '

# Row 36
$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").Value = '4.84'
$ws.Range("E36").Value = '  -2.59%  '

# Row 37
$ws.Range("D37").Value = '1.79'

# Row 38
$ws.Range("D38").Value = '343.84'
$ws.Range("E38").Value = '  +2.04%  '

# Row 39
$ws.Range("D39").Value = '0.953'
$ws.Range("E39").Value = '  -4.74%  '

# Row 40
$ws.Range("D40").Value = '6.21'
$ws.Range("E40").Value = '  -1.39%  '

# Row 41
$ws.Range("D41").Value = '4.03'
$ws.Range("E41").Value = '  -3.93%  '

# Row 42
$ws.Range("D42").Value = '38.40'
$ws.Range("E42").Value = '  -1.91%  '

# Row 43
$ws.Range("D43").Value = '21.30'
$ws.Range("E43").Value = '  -2.11%  '

# Row 44
$ws.Range("D44").Value = '20.69'
$ws.Range("E44").Value = '  -3.78%  '

# Row 45
$ws.Range("D45").Value = '0.0574'
$ws.Range("E45").Value = '  -2.53%  '

# Row 46
$ws.Range("D46").Value = '0.625'
$ws.Range("E46").Value = '  -1.01%  '

# Row 47
$ws.Range("B47").Value = 'FirstDigitalUSD'
$ws.Range("C47").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D47").Value = '0.998'
$ws.Range("E47").Value = '  -0.18%  '

# Row 48
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").Value = '0.0986'
$ws.Range("E48").Value = '  -3.61%  '

# Row 49
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").Value = '131.33'
$ws.Range("E49").Value = '  -2.66%  '

# Row 50
$ws.Range("B50").Value = 'WhiteBITCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D50").Value = '11.07'
$ws.Range("E50").Value = '  +0.31%  '

# Row 51
$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").Value = '0.0246'
$ws.Range("E51").Value = '  -3.88%  '
